$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 7457.846
$ws.Range("I76").Value = 10646.462
$ws.Range("J76").Value = 4269.231
$ws.Range("K76").Value = 10646.462
$ws.Range("L76").Value = 4269.231
$ws.Range("M76").Value = -10331.462
$ws.Range("N76").Value = -4899.231
$ws.Range("H79").Value = 7457.846
$ws.Range("I79").Value = 10646.462
$ws.Range("J79").Value = 4269.231
$ws.Range("K79").Value = 10646.462
$ws.Range("L79").Value = 4269.231
$ws.Range("M79").Value = -9554.462
$ws.Range("N79").Value = -6453.231
$ws.Range("H112").Value = 2075.4827
$ws.Range("I112").Value = 879.0909
$ws.Range("J112").Value = 2806.611
$ws.Range("K112").Value = 2637.2727
$ws.Range("L112").Value = 8419.832999999999
$ws.Range("M112").Value = -1529.2727
$ws.Range("N112").Value = -10635.833
$ws.Range("H133").Value = 88000
$ws.Range("J133").Value = 88000
$ws.Range("L133").Value = 88000
$ws.Range("N133").Value = -98120
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3145.2222
$ws.Range("I61").Value = 1962.4546
$ws.Range("J61").Value = 5003.857
$ws.Range("K61").Value = 1962.4546
$ws.Range("L61").Value = 5003.857
$ws.Range("M61").Value = -1750.4546
$ws.Range("N61").Value = -5427.857
$ws.Range("H122").Value = 2282
$ws.Range("I122").Value = 1692.8889
$ws.Range("J122").Value = 3165.6667
$ws.Range("K122").Value = 5078.6667
$ws.Range("L122").Value = 9497.000100000001
$ws.Range("M122").Value = -2628.6667
$ws.Range("N122").Value = -14397.0001
$ws.Range("H136").Value = 3145.2222
$ws.Range("I136").Value = 1962.4546
$ws.Range("J136").Value = 5003.857
$ws.Range("K136").Value = 5887.3638
$ws.Range("L136").Value = 15011.571
$ws.Range("M136").Value = -3337.3638
$ws.Range("N136").Value = -20111.571
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 8367.223
$ws.Range("I8").Value = 2375
$ws.Range("J8").Value = 13161
$ws.Range("K8").Value = 2375
$ws.Range("L8").Value = 13161
$ws.Range("M8").Value = -2235
$ws.Range("N8").Value = -13441
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = 0
$ws.Range("H15").Value = 8000
$ws.Range("J15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("N15").Value = -8454
$ws.Range("H134").Value = 2297.6667
$ws.Range("I134").Value = 2362.55
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 7087.650000000001
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -4552.650000000001
$ws.Range("N134").Value = -8070
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2394.1538
$ws.Range("I31").Value = 1883.1852
$ws.Range("J31").Value = 2757.2104
$ws.Range("K31").Value = 1883.1852
$ws.Range("L31").Value = 2757.2104
$ws.Range("M31").Value = -1588.1852
$ws.Range("N31").Value = -3347.2104
$ws.Range("H34").Value = 2394.1538
$ws.Range("I34").Value = 1883.1852
$ws.Range("J34").Value = 2757.2104
$ws.Range("K34").Value = 1883.1852
$ws.Range("L34").Value = 2757.2104
$ws.Range("M34").Value = -1681.1852
$ws.Range("N34").Value = -3161.2104
$ws.Range("H58").Value = 3324.4
$ws.Range("I58").Value = 3815.2
$ws.Range("J58").Value = 2342.8
$ws.Range("K58").Value = 3815.2
$ws.Range("L58").Value = 2342.8
$ws.Range("M58").Value = -3612.2
$ws.Range("N58").Value = -2748.8
$ws.Range("H105").Value = 4900.6313
$ws.Range("I105").Value = 4660.1
$ws.Range("J105").Value = 5167.8887
$ws.Range("K105").Value = 4660.1
$ws.Range("L105").Value = 5167.8887
$ws.Range("M105").Value = -2913.1
$ws.Range("N105").Value = -8661.8887
$ws.Range("H136").Value = 3324.4
$ws.Range("I136").Value = 3815.2
$ws.Range("J136").Value = 2342.8
$ws.Range("K136").Value = 11445.6
$ws.Range("L136").Value = 7028.400000000001
$ws.Range("M136").Value = -8895.599999999999
$ws.Range("N136").Value = -12128.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2242.814
$ws.Range("I68").Value = 1951.5454
$ws.Range("J68").Value = 2547.9524
$ws.Range("K68").Value = 5854.6362
$ws.Range("L68").Value = 7643.8572
$ws.Range("M68").Value = -5043.6362
$ws.Range("N68").Value = -9265.8572
$ws.Range("H69").Value = 4028.182
$ws.Range("I69").Value = 628
$ws.Range("K69").Value = 1884
$ws.Range("M69").Value = -1073
$ws.Range("H71").Value = 2242.814
$ws.Range("I71").Value = 1951.5454
$ws.Range("J71").Value = 2547.9524
$ws.Range("K71").Value = 17563.9086
$ws.Range("L71").Value = 22931.5716
$ws.Range("M71").Value = -13507.9086
$ws.Range("N71").Value = -31043.5716
$ws.Range("H72").Value = 4028.182
$ws.Range("I72").Value = 628
$ws.Range("K72").Value = 5652
$ws.Range("M72").Value = -1596
$ws.Range("H86").Value = 1973
$ws.Range("I86").Value = 700.5
$ws.Range("J86").Value = 2255.7778
$ws.Range("K86").Value = 2101.5
$ws.Range("L86").Value = 6767.3334
$ws.Range("M86").Value = -915.5
$ws.Range("N86").Value = -9139.3334
$ws.Range("H89").Value = 1973
$ws.Range("I89").Value = 700.5
$ws.Range("J89").Value = 2255.7778
$ws.Range("K89").Value = 6304.5
$ws.Range("L89").Value = 20302.0002
$ws.Range("M89").Value = -376.5
$ws.Range("N89").Value = -32158.0002
$ws.Range("H131").Value = 948.37036
$ws.Range("I131").Value = 437.5
$ws.Range("J131").Value = 1037.2174
$ws.Range("K131").Value = 1312.5
$ws.Range("L131").Value = 3111.6522
$ws.Range("M131").Value = 3727.5
$ws.Range("N131").Value = -13191.6522
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27519.234
$ws.Range("I70").Value = 47415.8
$ws.Range("J70").Value = 4909.5
$ws.Range("K70").Value = 47415.8
$ws.Range("L70").Value = 4909.5
$ws.Range("M70").Value = -47145.8
$ws.Range("N70").Value = -5449.5
$ws.Range("H73").Value = 27519.234
$ws.Range("I73").Value = 47415.8
$ws.Range("J73").Value = 4909.5
$ws.Range("K73").Value = 47415.8
$ws.Range("L73").Value = 4909.5
$ws.Range("M73").Value = -46479.8
$ws.Range("N73").Value = -6781.5
$ws.Range("H122").Value = 3008.64
$ws.Range("I122").Value = 2528.9092
$ws.Range("J122").Value = 6526.6665
$ws.Range("K122").Value = 7586.7276
$ws.Range("L122").Value = 19579.9995
$ws.Range("M122").Value = -5136.7276
$ws.Range("N122").Value = -24479.9995
$ws.Range("H132").Value = 2227
$ws.Range("I132").Value = 1739.5769
$ws.Range("J132").Value = 4037.4285
$ws.Range("K132").Value = 5218.7307
$ws.Range("L132").Value = 12112.2855
$ws.Range("M132").Value = -2688.7307
$ws.Range("N132").Value = -17172.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1862.7142
$ws.Range("I132").Value = 1477.5652
$ws.Range("J132").Value = 2600.9167
$ws.Range("K132").Value = 4432.6956
$ws.Range("L132").Value = 7802.750100000001
$ws.Range("M132").Value = -1902.6956
$ws.Range("N132").Value = -12862.7501
$ws.Range("H136").Value = 1609.0435
$ws.Range("I136").Value = 1506.2778
$ws.Range("J136").Value = 1979
$ws.Range("K136").Value = 4518.8334
$ws.Range("L136").Value = 5937
$ws.Range("M136").Value = -1968.8334
$ws.Range("N136").Value = -11037
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2689.1924
$ws.Range("I132").Value = 2651.375
$ws.Range("J132").Value = 2815.25
$ws.Range("K132").Value = 7954.125
$ws.Range("L132").Value = 8445.75
$ws.Range("M132").Value = -5424.125
$ws.Range("N132").Value = -13505.75
$ws.Range("H136").Value = 1567.7273
$ws.Range("I136").Value = 1302.5
$ws.Range("J136").Value = 3247.5
$ws.Range("K136").Value = 3907.5
$ws.Range("L136").Value = 9742.5
$ws.Range("M136").Value = -1357.5
$ws.Range("N136").Value = -14842.5
